$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.155.01"
$ws.Range("E2").Value = "  -3.64%  "
$ws.Range("D3").Value = "2.204.08"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "106.41"
$ws.Range("E5").Value = "  -14.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.29"
$ws.Range("E6").Value = "  +11.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  -3.06%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").Value = "  -6.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.22"
$ws.Range("E10").Value = "  -10.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0904"
$ws.Range("E11").Value = "  -4.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.15"
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.67"
$ws.Range("E13").Value = "  -7.70%  "
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.931"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.81"
$ws.Range("E16").Value = "  -4.39%  "
$ws.Range("D17").Value = "2.533.18"
$ws.Range("E17").Value = "  -3.47%  "
$ws.Range("D18").Value = "2.221.51"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "41.883.80"
$ws.Range("E19").Value = "  -4.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.27"
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("E21").Value = "  -5.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.82"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  +21.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.25"
$ws.Range("E24").Value = "  -7.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "226.43"
$ws.Range("E25").Value = "  -3.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.84"
$ws.Range("E26").Value = "  -7.20%  "
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.42"
$ws.Range("E28").Value = "  -4.86%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.66"
$ws.Range("E31").Value = "  -10.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.19"
$ws.Range("E32").Value = "  -5.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "172.50"
$ws.Range("E33").Value = "  -0.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.73"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0869"
$ws.Range("E35").Value = "  -5.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.46"
$ws.Range("E36").Value = "  -5.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.96"
$ws.Range("E37").Value = "  +6.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.27"
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.125"
$ws.Range("E39").Value = "  -4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0358"
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("E41").Value = "  -5.42%  "
$ws.Range("E42").Value = "  -4.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "69.92"
$ws.Range("E43").Value = "  -5.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.227"
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.42"
$ws.Range("E46").Value = "  -11.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.29"
$ws.Range("E47").Value = "  -6.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.36"
$ws.Range("E48").Value = "  -5.34%  "
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.66"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.38"
$ws.Range("E51").Value = "  -2.40%  "
